$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed by Excel as a
# number (single-dot decimal strings) need to be forced to Text format first
# so the exact original string (including trailing zeros) is preserved.

$ws.Range("D2").Value = "30.334.42"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").Value = "1.870.49"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.54"
$ws.Range("E5").Value = "  +1.15%  "

$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4688"
$ws.Range("E7").Value = "  +0.65%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2845"
$ws.Range("E8").Value = "  +1.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06544"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.49"
$ws.Range("E10").Value = "  +7.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07893"
$ws.Range("E11").Value = "  +1.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.10"
$ws.Range("E12").Value = "  +2.06%  "

$ws.Range("D13").Value = "1.870.72"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("E14").Value = "  +1.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6761"
$ws.Range("E15").Value = "  +1.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.18"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "30.327.78"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("E19").Value = "  +2.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.474"
$ws.Range("E20").Value = "  +3.14%  "

$ws.Range("D21").Value = "2.116.32"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007318"
$ws.Range("E22").Value = "  +1.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.152"
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.35"
$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.173"
$ws.Range("E26").Value = "  -0.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.21"
$ws.Range("E27").Value = "  +1.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.931"
$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09705"
$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.402"
$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("E32").Value = "  +0.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.105"
$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04717"
$ws.Range("E34").Value = "  +1.68%  "

$ws.Range("E35").Value = "  +4.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7067"
$ws.Range("E36").Value = "  +1.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.728"
$ws.Range("E37").Value = "  +0.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01865"
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.333"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.540"
$ws.Range("E40").Value = "  +1.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.70"
$ws.Range("E41").Value = "  +5.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.948"
$ws.Range("E42").Value = "  +1.05%  "

$ws.Range("E43").Value = "  -0.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4192"
$ws.Range("E44").Value = "  +1.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.72"
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.223"
$ws.Range("E47").Value = "  +0.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.261"
$ws.Range("E48").Value = "  +2.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "939.13"
$ws.Range("E49").Value = "  -3.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.28"
$ws.Range("E50").Value = "  +1.67%  "

$ws.Range("E51").Value = "  -0.56%  "
